# Apply updated crypto price/volume figures to Sheet1 (columns D and E).
# Values are written as text (matching the source inlineStr cells) by
# temporarily forcing a text number format, then restoring the default
# "Normal" style so no stray formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "65.217.41"
Set-TextValue $ws.Range("E2") "  -0.47%  "
Set-TextValue $ws.Range("D3") "3.562.30"
Set-TextValue $ws.Range("E3") "  -0.15%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.18%  "
Set-TextValue $ws.Range("D5") "600.17"
Set-TextValue $ws.Range("E5") "  +0.16%  "
Set-TextValue $ws.Range("D6") "133.84"
Set-TextValue $ws.Range("E6") "  -4.80%  "
Set-TextValue $ws.Range("D7") "3.559.08"
Set-TextValue $ws.Range("E7") "  -0.21%  "
Set-TextValue $ws.Range("E8") "  -0.05%  "
Set-TextValue $ws.Range("E9") "  -0.15%  "
Set-TextValue $ws.Range("E10") "  -2.13%  "
Set-TextValue $ws.Range("D11") "7.14"
Set-TextValue $ws.Range("E11") "  -0.07%  "
Set-TextValue $ws.Range("D12") "0.390"
Set-TextValue $ws.Range("E12") "  -0.90%  "
Set-TextValue $ws.Range("D13") "4.157.52"
Set-TextValue $ws.Range("E13") "  -0.43%  "
Set-TextValue $ws.Range("D14") "0.0000184"
Set-TextValue $ws.Range("E14") "  -2.66%  "
Set-TextValue $ws.Range("D15") "27.02"
Set-TextValue $ws.Range("E15") "  -0.46%  "
Set-TextValue $ws.Range("D16") "3.552.22"
Set-TextValue $ws.Range("E16") "  -0.44%  "
Set-TextValue $ws.Range("E17") "  -0.06%  "
Set-TextValue $ws.Range("D18") "65.273.88"
Set-TextValue $ws.Range("E18") "  -0.30%  "
Set-TextValue $ws.Range("D19") "9.94"
Set-TextValue $ws.Range("E19") "  -3.29%  "
Set-TextValue $ws.Range("D20") "14.45"
Set-TextValue $ws.Range("E20") "  +0.98%  "
Set-TextValue $ws.Range("D21") "5.86"
Set-TextValue $ws.Range("E21") "  -0.16%  "
Set-TextValue $ws.Range("D22") "392.64"
Set-TextValue $ws.Range("E22") "  -1.07%  "
Set-TextValue $ws.Range("D23") "0.579"
Set-TextValue $ws.Range("E23") "  +1.25%  "
Set-TextValue $ws.Range("D24") "3.696.98"
Set-TextValue $ws.Range("E24") "  -0.42%  "
Set-TextValue $ws.Range("D25") "74.16"
Set-TextValue $ws.Range("E26") "  +0.01%  "
Set-TextValue $ws.Range("E27") "  -2.04%  "
Set-TextValue $ws.Range("D28") "7.84"
Set-TextValue $ws.Range("E28") "  -0.52%  "
Set-TextValue $ws.Range("D29") "1.59"
Set-TextValue $ws.Range("E29") "  +25.99%  "
Set-TextValue $ws.Range("D30") "8.59"
Set-TextValue $ws.Range("E30") "  +3.68%  "
Set-TextValue $ws.Range("D31") "1.00"
Set-TextValue $ws.Range("E31") "  +0.06%  "
Set-TextValue $ws.Range("D32") "2.29"
Set-TextValue $ws.Range("E32") "  +0.44%  "
Set-TextValue $ws.Range("D33") "3.553.09"
Set-TextValue $ws.Range("E33") "  -0.98%  "
Set-TextValue $ws.Range("D34") "24.15"
Set-TextValue $ws.Range("E34") "  +0.89%  "
Set-TextValue $ws.Range("E35") "  -0.01%  "
Set-TextValue $ws.Range("D36") "0.146"
Set-TextValue $ws.Range("E36") "  -0.63%  "
Set-TextValue $ws.Range("D37") "170.73"
Set-TextValue $ws.Range("E37") "  +1.34%  "
Set-TextValue $ws.Range("D38") "6.96"
Set-TextValue $ws.Range("E38") "  -1.58%  "
Set-TextValue $ws.Range("D39") "1.55"
Set-TextValue $ws.Range("E39") "  -0.78%  "
Set-TextValue $ws.Range("D40") "5.08"
Set-TextValue $ws.Range("E40") "  +1.24%  "
Set-TextValue $ws.Range("D41") "0.0817"
Set-TextValue $ws.Range("E41") "  +1.58%  "
Set-TextValue $ws.Range("D42") "0.827"
Set-TextValue $ws.Range("E42") "  -1.74%  "
Set-TextValue $ws.Range("D43") "26.48"
Set-TextValue $ws.Range("E43") "  -0.64%  "
Set-TextValue $ws.Range("D44") "1.26"
Set-TextValue $ws.Range("E44") "  +5.37%  "
Set-TextValue $ws.Range("D45") "43.06"
Set-TextValue $ws.Range("E45") "  +0.09%  "
Set-TextValue $ws.Range("D46") "1.00"
Set-TextValue $ws.Range("E46") "  -0.08%  "
Set-TextValue $ws.Range("D47") "4.46"
Set-TextValue $ws.Range("E47") "  +0.27%  "
Set-TextValue $ws.Range("D48") "1.66"
Set-TextValue $ws.Range("E48") "  -2.54%  "
Set-TextValue $ws.Range("D49") "6.94"
Set-TextValue $ws.Range("E49") "  +1.59%  "
Set-TextValue $ws.Range("D50") "2.441.34"
Set-TextValue $ws.Range("E50") "  -0.20%  "
Set-TextValue $ws.Range("E51") "  -0.28%  "
